$wb = $excel.ActiveWorkbook

# --- Sheet1: "gramatge y espesor" -------------------------------------------------
$ws1 = $wb.Worksheets.Item("gramatge y espesor")

# Insert a new column B ("Gramage") before the existing espesor columns.
$ws1.Columns.Item(2).Insert()

# Header row
$ws1.Range("A1").Value = "Masa"
$ws1.Range("B1").Value = "Gramage"
$ws1.Range("C1").Value = "espesor 1"
$ws1.Range("D1").Value = "espesor 2"
$ws1.Range("E1").Value = "espesor 3"

# New Gramage column formula = Masa / 0.01
$ws1.Range("B2:B6").Formula = "=A2/0.01"

# --- Sheet2: "Espesor Total" -------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Espesor Total")

$ws2.Range("A1").Value = "espesor de paquete"
$ws2.Range("B1").Value = "Espesor individual de paquete"

$ws2.Range("B2:B4").Formula = "=A2/5"

# --- Sheet3: "Densidad" (new) -------------------------------------------------------
$ws3 = $wb.Worksheets.Add()
$ws3.Name = "Densidad"

$ws3.Range("A1").Value = "Densidad individual"
$ws3.Range("B1").Value = "Densidad individal de paquete"

$ws3.Range("A2").Formula = "='gramatge y espesor'!B2/AVERAGE('gramatge y espesor'!C2:E2)"
$ws3.Range("A3").Formula = "='gramatge y espesor'!B3/AVERAGE('gramatge y espesor'!C3:E3)"
$ws3.Range("A4").Formula = "='gramatge y espesor'!B4/AVERAGE('gramatge y espesor'!C4:E4)"
$ws3.Range("A5").Formula = "='gramatge y espesor'!B5/AVERAGE('gramatge y espesor'!C5:E5)"

$ws3.Range("B2").Formula = "='gramatge y espesor'!B2/AVERAGE('Espesor Total'!`$B`$2:`$B`$4)"
$ws3.Range("B3").Formula = "='gramatge y espesor'!B3/AVERAGE('Espesor Total'!`$B`$2:`$B`$4)"
$ws3.Range("B4").Formula = "='gramatge y espesor'!B4/AVERAGE('Espesor Total'!`$B`$2:`$B`$4)"
$ws3.Range("B5").Formula = "='gramatge y espesor'!B5/AVERAGE('Espesor Total'!`$B`$2:`$B`$4)"

# Move the new sheet to the end (after "Espesor Total") to match expected sheet order.
# (Re-fetch the Espesor Total reference AFTER the new sheet was added.)
$wsEspesorTotal = $wb.Worksheets.Item("Espesor Total")
$ws3.Move($null, $wsEspesorTotal)

$wb.Save()
